$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the "Title" column (column A) entirely, shifting B:D left to A:C
$ws.Range("A1").EntireColumn.Delete()

# Update the selection to match the post-edit state (selecting column A)
$ws.Range("A1:A1048576").Select()
